$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update name from "Anu" to "Sandeep"
$ws.Range("A2").Value = "Sandeep"

# Update password / confirm password from "Anu@1234" to "Sandeep@123"
$ws.Range("F2").Value = "Sandeep@123"
$ws.Range("G2").Value = "Sandeep@123"

# Update phone number (input time error fix) from 9039568978 to 7817001234
$ws.Range("C2").Value = 7817001234

# Move active selection to H2
$ws.Range("H2").Select()
